# Weekly update: insert 3 new price records (week of 45006) at the top of the
# "Membrillo" price table, pushing the existing rows down by three positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at the former first data-row position (row 108).
# Everything that used to live in rows 108:136 now lives in rows 111:139.
$ws.Rows("108:110").Insert()

# New weekly records to populate in the freshly inserted rows.
$newRows = @(
  @{ Row = 108; L = "Especial"; M = 12; N = 200000; O = 200000; P = 200000; S = 444 },
  @{ Row = 109; L = "Primera";  M = 20; N = 180000; O = 180000; P = 180000; S = 400 },
  @{ Row = 110; L = "Segunda";  M = 16; N = 150000; O = 150000; P = 150000; S = 333 }
)

foreach ($r in $newRows) {
  $row = $r.Row

  $ws.Range("A$row").Value = 6
  $ws.Range("B$row").Value = "Mercado Mayorista Lo Valledor de Santiago"
  $ws.Range("C$row").Value = "Metropolitana"
  $ws.Range("D$row").Value = 45006
  $ws.Range("E$row").Value = 13
  $ws.Range("F$row").Value = "Fruta"
  $ws.Range("G$row").Value = 100104
  $ws.Range("H$row").Value = "Frutos de pepita"
  $ws.Range("I$row").Value = 100104003
  $ws.Range("J$row").Value = "Membrillo"
  $ws.Range("K$row").Value = "Champion"
  $ws.Range("L$row").Value = $r.L
  $ws.Range("M$row").Value = $r.M
  $ws.Range("N$row").Value = $r.N
  $ws.Range("O$row").Value = $r.O
  $ws.Range("P$row").Value = $r.P
  $ws.Range("Q$row").Value = "$/bins (450 kilos)"
  $ws.Range("R$row").Value = "Región de O'Higgins"
  $ws.Range("S$row").Value = $r.S
  $ws.Range("T$row").Value = 450
}
